$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New yellow separator row (row 26) matching the other section breaks ---
# (applied first so the plain "yellow fill, no border" style is registered
# before the bordered variant used below)
$ws.Range("D26:F26").Interior.Color = 65535

# --- New "Random button" assertion row (row 17, column B) ---
# Yellow highlight fill + full thin box border around the new assertion cell.
$ws.Range("B17").Interior.Color = 65535
$ws.Range("B17").Borders.LineStyle = 1
$ws.Range("B17").Borders.Weight = 2
$ws.Range("B17").Value = "Random button"

# --- New user story / conditional / manual test row (row 27) ---
$ws.Range("D27").Value = "As a user, the app can suggest a restaurant for me"
$ws.Range("E27").Value = 'After clicking the "Surprise Me" button, a random restaurant is selected and presented to me'
$ws.Range("F27").Value = "Clicked the button multiple times and it opened different restaurants."

# --- Update the view: scroll so column F is leftmost, select F28 ---
[void]$ws.Activate()
[void]$ws.Range("F28").Select()
$excel.ActiveWindow.ScrollColumn = 6
